{"js": "// Change 1: collapse \"What is the \" + \"current status\" (+ proofErr marks) +\n// \" of the microservice? Hopefully, it's done!\" into a single run with the\n// combined text (no grammar-check proof markers left behind).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet statusPara = null;\nfor (const para of paragraphs.items) {\n  if (para.text.indexOf(\"What is the\") !== -1 && para.text.indexOf(\"current status\") !== -1) {\n    statusPara = para;\n    break;\n  }\n}\n\nif (statusPara) {\n  statusPara.getRange().insertText(\n    \"What is the current status of the microservice? Hopefully, it\\u2019s done!\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// Change 2: split the \"I am assuming ... deployment. Other than that ...\"\n// run so that a new clause (\" likely utilizing a hosting service for the\n// microservice\") is inserted right after \"deployment\" and before the\n// following period.\nlet assumingPara = null;\nfor (const para of paragraphs.items) {\n  if (para.text.indexOf(\"I am assuming that my teammates\") !== -1) {\n    assumingPara = para;\n    break;\n  }\n}\n\nif (assumingPara) {\n  const deploymentMatches = assumingPara.search(\"deployment\", { matchCase: true });\n  deploymentMatches.load(\"text\");\n  await context.sync();\n\n  const deploymentRange = deploymentMatches.items[0];\n  const insertAfterDeployment = deploymentRange.getRange(\"End\");\n  const insertedRange = insertAfterDeployment.insertText(\n    \" likely utilizing a hosting service for the microservice\",\n    Word.InsertLocation.replace\n  );\n  // Re-apply the run's formatting explicitly so the inserted text becomes\n  // its own run (matching the paragraph's existing non-bold body style)\n  // instead of silently merging with its neighbour.\n  insertedRange.font.set({\n    name: \"Calibri\",\n    nameAscii: \"Calibri\",\n    nameOther: \"Calibri\",\n    nameBidirectional: \"Calibri\",\n    size: 12\n  });\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Change 1: collapse \"What is the \" + \"current status\" (+ proofErr marks) +\n# \" of the microservice? Hopefully, it's done!\" into a single clean run.\n# Find/Replace across the whole matched span merges the three original runs\n# (and drops the grammar-check proofErr bookends) into one replacement run.\n$apostrophe = [char]0x2019\n$statusText = \"What is the current status of the microservice? Hopefully, it\" + $apostrophe + \"s done!\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $statusText\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $statusText\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# Change 2: split the \"I am assuming ... deployment. Other than that ...\"\n# run so a new clause (\" likely utilizing a hosting service for the\n# microservice\") is inserted right after \"deployment\" and before the period\n# that follows it.\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $candidate = $d.Paragraphs.Item($i)\n    if ($candidate.Range.Text -like \"I am assuming that my teammates*\") {\n        $target = $candidate\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $scoped = $target.Range.Duplicate\n    $scoped.Find.Execute(\"deployment\") | Out-Null\n    $scoped.Collapse(0)\n    $scoped.InsertAfter(\" likely utilizing a hosting service for the microservice\")\n    # Re-apply the run's formatting explicitly so the inserted text lands in\n    # its own run (matching the paragraph's existing non-bold body style)\n    # rather than silently merging into its neighbour run.\n    $scoped.Font.Name = \"Calibri\"\n    $scoped.Font.NameAscii = \"Calibri\"\n    $scoped.Font.NameOther = \"Calibri\"\n    $scoped.Font.NameBi = \"Calibri\"\n    $scoped.Font.Size = 12\n}\n"}
